# This script reproduces the weekly data refresh for the
# "Brocoli - Vega Monumental Concepcion" sheet:
#   - Two new daily records (date serial 44813) are added for
#     "Primera" and "Segunda" quality, pushing the previously
#     existing data (which used to start at row 228) down by two
#     rows.
#   - Because the two newly inserted rows carry exactly the same
#     values that used to occupy rows 228/229 (only the date is
#     genuinely new), the easiest faithful reproduction is:
#       1) remember the current contents of rows 228 and 229
#       2) insert two blank rows right before the old row 230
#          (this naturally shifts every following row down by 2,
#          growing the sheet from 352 to 354 rows)
#       3) put the remembered old 228/229 values into the freshly
#          inserted rows 230/231
#       4) stamp the new date (44813) onto rows 228 and 229

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: capture the existing row 228 / 229 data (columns A:R)
$oldRow228 = $ws.Range("A228:R228").Value2
$oldRow229 = $ws.Range("A229:R229").Value2

# Step 2: insert two new rows right after row 229 (i.e. before the
# old row 230), shifting all subsequent rows down by two.
$ws.Rows("230:231").Insert()

# Step 3: fill the newly inserted rows with what used to be in 228/229
$ws.Range("A230:R230").Value2 = $oldRow228
$ws.Range("A231:R231").Value2 = $oldRow229

# Step 4: update the date on the first two rows to the new reading date
$ws.Range("D228").Value2 = 44813
$ws.Range("D229").Value2 = 44813
